$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the "<add>" and "</add>" marker runs that bracket the
#    commented "b" (the b/p correction note, comment id 0). These are
#    the very first occurrences of those tokens in the document.
# ------------------------------------------------------------------
$d.Content.Find.Execute("<add>", $true, $false, $false, $false, $false, $true, 1, $false, "", 1) | Out-Null
$d.Content.Find.Execute("</add>", $true, $false, $false, $false, $false, $true, 1, $false, "", 1) | Out-Null

# ------------------------------------------------------------------
# 2. Simple spelling / wording fixes (single w:t runs).
# ------------------------------------------------------------------
$d.Content.Find.Execute("encores avec l", $true, $false, $false, $false, $false, $true, 1, $false, "encores avecq l", 1) | Out-Null
$d.Content.Find.Execute("y estre se brusle avec luy qui ne se brusle point Cela", $true, $false, $false, $false, $false, $true, 1, $false, "y estre se brusle avecq luy qui ne se brusle point Cela", 1) | Out-Null
$d.Content.Find.Execute("invenction Pile le da", $true, $false, $false, $false, $false, $true, 1, $false, "invention Pile le da", 1) | Out-Null

# ------------------------------------------------------------------
# 3. "peu a chasque foye" -> "peu a chasque foy" + a separate trailing
#    run containing just "s" (its own, unstyled run rather than being
#    merged back into the preceding run).
# ------------------------------------------------------------------
$d.Content.Find.Execute("peu a chasque foye", $true, $false, $false, $false, $false, $true, 1, $false, "peu a chasque foy", 1) | Out-Null

$text = $d.Content.Text
$idx = $text.IndexOf("peu a chasque foy")
$endPos = $idx + "peu a chasque foy".Length

# Build the "s" text in a throwaway paragraph first so it is created
# with plain/default run formatting (no inherited color), then move
# it (cut/paste) into place right after "...foy" so it lands as its
# own run instead of merging into the neighbouring colored run.
$scratchPara = $d.Paragraphs.Add()
$scratchPara.Range.Text = "s"
$scratchRange = $d.Range($scratchPara.Range.Start, $scratchPara.Range.Start + 1)
$scratchRange.Cut()

$target = $d.Range($endPos, $endPos)
$target.Paste()

# Clean up the now-empty scratch paragraph (and its paragraph mark) we
# used as scaffolding, restoring the document to its original
# paragraph count.
$lastPara = $d.Paragraphs.Last
$cleanupRange = $d.Range($lastPara.Range.Start - 1, $lastPara.Range.End)
$cleanupRange.Delete()

Write-Output "edits applied"
